$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apache POI")
$ws.Range("D4").Value = "Apache POI"
